# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
#
# For each of the following row pairs, the two match-rows were swapped in
# their entirety (every column from B through AD - i.e. id, Div, Date,
# HomeTeam, AwayTeam, full-time/half-time goals, result, and all odds
# columns) while the leading rank column "A" stays put for each physical
# row. Column D (Date) is identical between each paired row anyway, so
# including it in the swap is a no-op.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(69, 70),
    @(78, 79),
    @(89, 90),
    @(103, 104),
    @(108, 109),
    @(223, 224),
    @(231, 232)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
